$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-18: update B, C, D columns with new values
$ws.Range("B2").Value = 65724
$ws.Range("C2").Value = 5632.2595
$ws.Range("D2").Value = 60091.7405

$ws.Range("B3").Value = 63085
$ws.Range("C3").Value = 5478.844000000001
$ws.Range("D3").Value = 57606.156

$ws.Range("B4").Value = 61258
$ws.Range("C4").Value = 5431.737499999999
$ws.Range("D4").Value = 55826.2625

$ws.Range("B5").Value = 60272
$ws.Range("C5").Value = 5351.2935
$ws.Range("D5").Value = 54920.7065

$ws.Range("B6").Value = 63182
$ws.Range("C6").Value = 5443.7565
$ws.Range("D6").Value = 57738.2435

$ws.Range("B7").Value = 67802
$ws.Range("C7").Value = 5521.645500000001
$ws.Range("D7").Value = 62280.3545

$ws.Range("B8").Value = 69571
$ws.Range("C8").Value = 6568.996000000001
$ws.Range("D8").Value = 63002.004

$ws.Range("B9").Value = 81158
$ws.Range("C9").Value = 7039.0705
$ws.Range("D9").Value = 74118.9295

$ws.Range("B10").Value = 98453
$ws.Range("C10").Value = 8768.984
$ws.Range("D10").Value = 89684.016

$ws.Range("B11").Value = 108429
$ws.Range("C11").Value = 14287.7385
$ws.Range("D11").Value = 94141.26149999999

$ws.Range("B12").Value = 111739
$ws.Range("C12").Value = 16203.32
$ws.Range("D12").Value = 95535.67999999999

$ws.Range("B13").Value = 113097
$ws.Range("C13").Value = 16253.0025
$ws.Range("D13").Value = 96843.9975

$ws.Range("B14").Value = 112752
$ws.Range("C14").Value = 16193.9575
$ws.Range("D14").Value = 96558.0425

$ws.Range("B15").Value = 116975
$ws.Range("C15").Value = 16307.025
$ws.Range("D15").Value = 100667.975

$ws.Range("B16").Value = 117642
$ws.Range("C16").Value = 16253.9055
$ws.Range("D16").Value = 101388.0945

$ws.Range("B17").Value = 96422
$ws.Range("C17").Value = 15667.071
$ws.Range("D17").Value = 80754.929

$ws.Range("B18").Value = 94648
$ws.Range("C18").Value = 16597.392
$ws.Range("D18").Value = 78050.60800000001

# Rows 19-25: delete B column cells entirely, update C column; D stays the same
$ws.Range("B19").ClearContents()
$ws.Range("C19").Value = 16121.5775

$ws.Range("B20").ClearContents()
$ws.Range("C20").Value = 15277.878

$ws.Range("B21").ClearContents()
$ws.Range("C21").Value = 13983.214

$ws.Range("B22").ClearContents()
$ws.Range("C22").Value = 12008.9865

$ws.Range("B23").ClearContents()
$ws.Range("C23").Value = 9648.002000000002

$ws.Range("B24").ClearContents()
$ws.Range("C24").Value = 6756.816500000001

$ws.Range("B25").ClearContents()
$ws.Range("C25").Value = 6005.93
